$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 192
$ws.Range("F5").Value = 369
$ws.Range("F9").Value = 4216
$ws.Range("F10").Value = 4216
$ws.Range("F14").Value = 6117
$ws.Range("F15").Value = 65
$ws.Range("F17").Value = 2344
$ws.Range("F20").Value = 478
$ws.Range("F21").Value = 9218
$ws.Range("F23").Value = 2478
$ws.Range("F25").Value = 2317
$ws.Range("F26").Value = 2457
$ws.Range("F28").Value = 243
$ws.Range("F29").Value = 1972
$ws.Range("F31").Value = 60
$ws.Range("F32").Value = 331
$ws.Range("F34").Value = 44
$ws.Range("F36").Value = 44
$ws.Range("F37").Value = 62
$ws.Range("F39").Value = 1222
$ws.Range("F40").Value = 1220
$ws.Range("F44").Value = 1544
$ws.Range("F45").Value = 2541
$ws.Range("F46").Value = 926
$ws.Range("F47").Value = 303
$ws.Range("F48").Value = 1253
$ws.Range("F49").Value = 24

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 150
$ws.Range("F22").Value = 75

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 694
$ws.Range("F3").Value = 898
$ws.Range("F4").Value = 102

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 192
$ws.Range("F3").Value = 694
$ws.Range("F4").Value = 898
$ws.Range("F6").Value = 369
$ws.Range("F14").Value = 4216
$ws.Range("F18").Value = 6117
$ws.Range("F19").Value = 65
$ws.Range("F22").Value = 478
$ws.Range("F23").Value = 9218
$ws.Range("F25").Value = 2478
$ws.Range("F27").Value = 2317
$ws.Range("F28").Value = 2457
$ws.Range("F30").Value = 243
$ws.Range("F31").Value = 1972
$ws.Range("F33").Value = 60
$ws.Range("F34").Value = 331
$ws.Range("F36").Value = 62
$ws.Range("F38").Value = 1220
$ws.Range("F42").Value = 1544
$ws.Range("F43").Value = 2541
$ws.Range("F44").Value = 926
$ws.Range("F45").Value = 303
$ws.Range("F48").Value = 1253
$ws.Range("F50").Value = 75
$ws.Range("F51").Value = 75
